# "remake tables and other predictions"
# Rebuild the prediction-accuracy table: add an "Old result" header row,
# split the existing year.* / random.* rows into two blocks (with the new
# "+quarterFE" / "+monthFE+quaterFE" variants inserted), leave two blank
# spacer rows between the blocks, and refresh layout/page-setup metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the old table body (rows 3:6) so we can rebuild it from scratch ---
$ws.Range("B3:E6").Clear()

# --- row 3: new, unstyled "Old result" summary row -------------------------
$ws.Range("B3").Value = "Old result "
$ws.Range("C3").Value = 0.536
$ws.Range("D3").Value = 0.623
$ws.Range("E3").Value = 0.169

# --- rows 5-8: year.* block --------------------------------------------------
$ws.Range("C5:E8").NumberFormat = "0.000"

$ws.Range("B5").Value = "year.ols"
$ws.Range("C5").Value = 0.5487049
$ws.Range("D5").Value = 0.5688843
$ws.Range("E5").Value = 0.1299255

$ws.Range("B6").Value = "year.ols+ quarterFE"
$ws.Range("C6").Value = 0.559955
$ws.Range("D6").Value = 0.6008912
$ws.Range("E6").Value = 0.1272171

$ws.Range("B7").Value = "year.ols+ monthFE+quaterFE"
$ws.Range("C7").Value = 0.56921231
$ws.Range("D7").Value = 0.5951778
$ws.Range("E7").Value = 0.1001644

$ws.Range("B8").Value = "year.LASSO"
$ws.Range("C8").Value = 0.56968424
$ws.Range("D8").Value = 0.6002007
$ws.Range("E8").Value = 0.11101

# --- rows 9-10: blank spacer rows, still carrying the numeric style ---------
$ws.Range("C9:E10").NumberFormat = "0.000"
$ws.Range("C9:E10").Value = ""

# --- rows 11-14: random.* block ---------------------------------------------
$ws.Range("C11:E14").NumberFormat = "0.000"

$ws.Range("B11").Value = "random.ols"
$ws.Range("C11").Value = 0.4970718
$ws.Range("D11").Value = 0.5477649
$ws.Range("E11").Value = 0.0950986

$ws.Range("B12").Value = "random.ols+ quarterFE"
$ws.Range("C12").Value = 0.5181456
$ws.Range("D12").Value = 0.5407152
$ws.Range("E12").Value = 0.1117934

$ws.Range("B13").Value = "random.ols+ monthFE+quaterFE"
$ws.Range("C13").Value = 0.526206
$ws.Range("D13").Value = 0.5419876
$ws.Range("E13").Value = 0.1802338

$ws.Range("B14").Value = "random.LASSO "
$ws.Range("C14").Value = 0.5078364
$ws.Range("D14").Value = 0.5592047
$ws.Range("E14").Value = 0.1662719

# --- layout: wider label column, page setup ---------------------------------
# (closest value the engine's char->pixel rounding can land on 31.625 units)
$ws.Columns("B").ColumnWidth = 30.833333333333332

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection parks below the rebuilt table, matching the saved file -------
[void]$ws.Range("E20").Select()
